$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")
$ws.Activate()

# Update min_hours_per_week (column I) for all employee rows from 44 to 36
$ws.Range("I2:I63").Value = 36

# Update consecutive_worked_sundays (column L) for rows 30-63 from 1 to 0
$ws.Range("L30:L63").Value = 0

$ws.Range("P23").Select()
